$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.6753301551942219, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 0, 2.020003488691574)
    3  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 5.553084769722144)
    4  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    5  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
    6  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797)
    7  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 3.781711156805759)
    8  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 3.645393585217082)
    9  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 3.645393585217082)
    10 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    11 = @(0.3048080303191223, 0.00007097389502863649, 0.1575252929769615, 0.496779210170732, 0, 0.9591835073618444)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
